# Junction_Flooding_461: replace the 1000-row sample window with a new
# window of readings ("custom accuracy + 데이터 1000개"), drop the now-unused
# 6th sample row, and widen a handful of data columns by one character.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-5, columns A (Time) .. AH (J33), replacing the
# existing sample readings in place.
$data = @(
    @(45142.50694444445, 14.315, 9.452, 3.548, 30.929, 23.407, 11.081, 33.367, 17.452, 7.022, 10.47, 12.134, 12.887, 3.618, 11.279, 15.513, 9.993, 3.077, 1.709, 164.3, 31.24, 10.411, 20.071, 10.298, 2.945, 17.697, 9.196, 8.398, 9.994999999999999, 12.384, 3.078, 30.175, 5.636, 13.015),
    @(45142.51388888889, 5.68, 3.717, 1.401, 12.546, 9.167, 4.35, 19.663, 6.981, 2.813, 3.925, 4.929, 5.29, 1.455, 4.512, 6.205, 4.233, 1.367, 0.707, 61.372, 12.799, 4.164, 8.103999999999999, 4.132, 1.23, 9.573, 3.678, 3.478, 4.108, 4.933, 1.199, 18.551, 2.167, 5.208),
    @(45142.52083333334, 21.063, 15.498, 1.459, 46.026, 37.207, 16.489, 60.523, 25.596, 11.266, 16.628, 18.41, 19.556, 5.312, 16.542, 23.446, 14.069, 1.039, 0.951, 244.464, 46.169, 15.269, 30.877, 16.199, 2.589, 30.073, 13.487, 12.013, 14.138, 19.193, 0.747, 55.031, 8.551, 19.09),
    @(45142.52777777778, 15.31, 11.28, 1.06, 33.48, 27.03, 11.98, 48.45, 18.62, 8.210000000000001, 12.06, 13.4, 14.25, 3.87, 12.03, 17.07, 10.25, 0.77, 0.68, 175.81, 33.7, 11.1, 22.52, 11.79, 1.89, 23.37, 9.81, 8.75, 10.29, 13.97, 0.54, 44.25, 6.21, 13.88)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}

# The old 6th sample row is no longer part of the window; removing it also
# shrinks the sheet dimension from AH6 down to AH5.
$ws.Rows.Item(6).Delete()

# A subset of the data columns become one character wider (stored width
# 7 -> 8). ColumnWidth is specified in characters and Excel adds ~5/6 of a
# character of internal padding when persisting the stored <col> width, so
# subtract that padding here to land exactly on a stored width of 8.
$wideColumns = @(3, 7, 10, 11, 17, 22, 24, 27, 28, 29)
foreach ($col in $wideColumns) {
    $ws.Columns.Item($col).ColumnWidth = 7.166666666666667
}
